# Insert a new data row at row 91 (pushes existing rows 91:178 down to 92:179),
# then populate the new row with a copy of the (now shifted-down) former row 91
# data, updated with the new week's Fecha (Date) and Volumen values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 91 and below down by one row.
$ws.Rows("91:91").Insert()

# Populate the newly inserted row 91.
$ws.Range("A91").Value = 8
$ws.Range("B91").Value = "Terminal La Palmera de La Serena"
$ws.Range("C91").Value = "Coquimbo"
$ws.Range("D91").Value = 45167
$ws.Range("E91").Value = 4
$ws.Range("F91").Value = 100114007
$ws.Range("G91").Value = "Jengibre"
$ws.Range("H91").Value = "Sin especificar"
$ws.Range("I91").Value = "Primera"
$ws.Range("J91").Value = 400
$ws.Range("K91").Value = 17000
$ws.Range("L91").Value = 18000
$ws.Range("M91").Value = 17500
$ws.Range("N91").Value = "$/caja 13 kilos"
$ws.Range("O91").Value = "Perú"
$ws.Range("P91").Value = 1346
$ws.Range("Q91").Value = 13
$ws.Range("R91").Value = "Hortaliza"
